$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values (Price / Volume columns) are stored as text,
# matching the source workbook which keeps these as literal strings (e.g. "289.21", "0.93%").
$textCells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "E21", "D23", "E23", "D24", "E24", "E25", "D26", "E26", "D27", "E27", "D28", "E28", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '290.59'
$ws.Range('E2').Value = '1.40%'
$ws.Range('D3').Value = '29.52'
$ws.Range('E3').Value = '4.15%'
$ws.Range('D4').Value = '5.104'
$ws.Range('E4').Value = '3.75%'
$ws.Range('D5').Value = '0.06688'
$ws.Range('E5').Value = '2.13%'
$ws.Range('D6').Value = '7.355'
$ws.Range('E6').Value = '1.50%'
$ws.Range('B7').Value = 'GateToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D7').Value = '3.406'
$ws.Range('E7').Value = '0.28%'
$ws.Range('B8').Value = 'FTXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D8').Value = '1.352'
$ws.Range('E8').Value = '-0.93%'
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D9').Value = '0.9160'
$ws.Range('E9').Value = '0.07%'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').Value = '0.1589'
$ws.Range('E10').Value = '1.21%'
$ws.Range('B11').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C11').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D11').Value = '0.06689'
$ws.Range('E11').Value = '0.47%'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').Value = '0.07680'
$ws.Range('E12').Value = '-0.11%'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').Value = '0.02941'
$ws.Range('E13').Value = '-1.30%'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').Value = '0.08992'
$ws.Range('E14').Value = '0.13%'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').Value = '0.001587'
$ws.Range('E15').Value = '-1.02%'
$ws.Range('B16').Value = 'CoinExToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D16').Value = '0.04518'
$ws.Range('E16').Value = '1.29%'
$ws.Range('B17').Value = 'One'
$ws.Range('C17').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D17').Value = '0.0006476'
$ws.Range('E17').Value = '-1.34%'
$ws.Range('B18').Value = 'TigerCash'
$ws.Range('C18').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D18').Value = '0.006280'
$ws.Range('E18').Value = '4.21%'
$ws.Range('B19').Value = 'LEO'
$ws.Range('C19').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D19').Value = '3.451'
$ws.Range('E19').Value = '-0.97%'
$ws.Range('D20').Value = '2.221'
$ws.Range('E20').Value = '-0.93%'
$ws.Range('E21').Value = '1.79%'
$ws.Range('D23').Value = '4.070'
$ws.Range('E23').Value = '2.33%'
$ws.Range('D24').Value = '0.1551'
$ws.Range('E24').Value = '2.05%'
$ws.Range('E25').Value = '0.39%'
$ws.Range('D26').Value = '0.004136'
$ws.Range('E26').Value = '-4.68%'
$ws.Range('D27').Value = '0.0001249'
$ws.Range('E27').Value = '5.86%'
$ws.Range('D28').Value = '0.0001617'
$ws.Range('E28').Value = '-1.11%'
$ws.Range('D40').Value = '0.04232'
$ws.Range('E40').Value = '2.03%'
$ws.Range('D41').Value = '0.006742'
$ws.Range('E41').Value = '0.35%'
$ws.Range('D42').Value = '0.1242'
$ws.Range('E42').Value = '-12.06%'
$ws.Range('D43').Value = '0.001978'
$ws.Range('E43').Value = '-3.93%'
$ws.Range('D44').Value = '0.01166'
$ws.Range('E44').Value = '-6.26%'
$ws.Range('D45').Value = '0.00005612'
$ws.Range('E45').Value = '1.02%'
$ws.Range('B46').Value = 'CoinbaseStockToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range('D46').Value = '0.01306'
$ws.Range('E46').Value = '-29.41%'
$ws.Range('B47').Value = 'BOLO'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range('D47').Value = '1.974'
$ws.Range('E47').Value = '26.46%'
